$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new price text is purely numeric-looking (e.g. "609.53") need their
# number format forced to Text first, otherwise Excel auto-converts the assigned
# string into a numeric value instead of keeping it as text (matching the source
# data, which stores these as text / inline strings).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "70.207.50"
$ws.Range("E2").Value = "  -0.20%  "
$ws.Range("D3").Value = "3.562.79"
$ws.Range("E3").Value = "  -0.26%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "609.53"
$ws.Range("E5").Value = "  +3.44%  "
$ws.Range("D6").Value = "185.95"
$ws.Range("E6").Value = "  -0.47%  "
$ws.Range("D7").Value = "3.555.43"
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("E8").Value = "  -0.33%  "
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("D10").Value = "0.215"
$ws.Range("E10").Value = "  +7.36%  "
$ws.Range("E11").Value = "  -0.63%  "
$ws.Range("D12").Value = "53.78"
$ws.Range("E12").Value = "  -1.90%  "
$ws.Range("D13").Value = "0.0000309"
$ws.Range("E13").Value = "  +0.42%  "
$ws.Range("D14").Value = "9.47"
$ws.Range("E14").Value = "  -0.69%  "
$ws.Range("D15").Value = "4.129.70"
$ws.Range("E15").Value = "  -0.25%  "
$ws.Range("D16").Value = "70.262.51"
$ws.Range("E16").Value = "  -0.09%  "
$ws.Range("D17").Value = "3.560.20"
$ws.Range("E17").Value = "  +0.05%  "
$ws.Range("D18").Value = "12.67"
$ws.Range("E18").Value = "  +1.25%  "
$ws.Range("D19").Value = "18.96"
$ws.Range("E19").Value = "  -2.52%  "
$ws.Range("D20").Value = "579.99"
$ws.Range("E20").Value = "  +6.51%  "
$ws.Range("E21").Value = "  +0.33%  "
$ws.Range("D22").Value = "0.993"
$ws.Range("E22").Value = "  -2.71%  "
$ws.Range("D23").Value = "17.35"
$ws.Range("E23").Value = "  -3.49%  "
$ws.Range("D24").Value = "4.71"
$ws.Range("E24").Value = "  +0.14%  "
$ws.Range("D25").Value = "4.86"
$ws.Range("E25").Value = "  -1.49%  "
$ws.Range("D26").Value = "94.83"
$ws.Range("E26").Value = "  -0.99%  "
$ws.Range("D27").Value = "2.95"
$ws.Range("E27").Value = "  -2.17%  "
$ws.Range("E28").Value = "  -4.51%  "
$ws.Range("D29").Value = "9.37"
$ws.Range("E29").Value = "  +2.49%  "
$ws.Range("D30").Value = "32.24"
$ws.Range("E30").Value = "  -0.14%  "
$ws.Range("D31").Value = "7.03"
$ws.Range("E31").Value = "  -4.25%  "
$ws.Range("D32").Value = "12.24"
$ws.Range("E32").Value = "  -2.53%  "
$ws.Range("E33").Value = "  -0.71%  "
$ws.Range("D34").Value = "63.52"
$ws.Range("E34").Value = "  -2.64%  "
$ws.Range("D35").Value = "3.68"
$ws.Range("E35").Value = "  +19.20%  "
$ws.Range("E36").Value = "  +0.00%  "
$ws.Range("D37").Value = "531.89"
$ws.Range("E37").Value = "  -3.69%  "
$ws.Range("E38").Value = "  -3.61%  "
$ws.Range("E39").Value = "  +0.18%  "
$ws.Range("D40").Value = "37.11"
$ws.Range("E40").Value = "  -3.55%  "
$ws.Range("D41").Value = "0.0₃0783"
$ws.Range("E41").Value = "  +1.95%  "
$ws.Range("D42").Value = "3.535.69"
$ws.Range("E42").Value = "  +4.91%  "
$ws.Range("D43").Value = "3.52"
$ws.Range("E43").Value = "  +4.22%  "
$ws.Range("E44").Value = "  +0.48%  "
$ws.Range("E45").Value = "  +2.96%  "
$ws.Range("D46").Value = "3.47"
$ws.Range("E46").Value = "  -2.85%  "
$ws.Range("D47").Value = "2.92"
$ws.Range("E47").Value = "  -1.98%  "
$ws.Range("E48").Value = "  +2.80%  "
$ws.Range("D49").Value = "9.21"
$ws.Range("E49").Value = "  -0.10%  "
$ws.Range("E50").Value = "  +0.25%  "
$ws.Range("D51").Value = "136.43"
$ws.Range("E51").Value = "  -0.75%  "

